$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting existing rows 87..201 down to 88..202
$ws.Rows("87").Insert()

# Populate the newly inserted row 87 with its data
$ws.Range("A87").Value = 3
$ws.Range("B87").Value = "Femacal de La Calera"
$ws.Range("C87").Value = "Coquimbo"
$ws.Range("D87").Value = 44467
$ws.Range("E87").Value = 5
$ws.Range("F87").Value = 100112043
$ws.Range("G87").Value = "Pepino ensalada"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 115
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 12500
$ws.Range("M87").Value = 12261
$ws.Range("N87").Value = "$/caja 70 unidades"
$ws.Range("O87").Value = "Región de Arica y Parinacota"
$ws.Range("P87").Value = 175
$ws.Range("Q87").Value = 70
$ws.Range("R87").Value = "Hortaliza"
